$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1767.4166
$ws.Range("I80").Value = 3556
$ws.Range("J80").Value = 489.85715
$ws.Range("K80").Value = 10668
$ws.Range("L80").Value = 1469.57145
$ws.Range("M80").Value = -9670
$ws.Range("N80").Value = -3465.57145
$ws.Range("H83").Value = 1767.4166
$ws.Range("I83").Value = 3556
$ws.Range("J83").Value = 489.85715
$ws.Range("K83").Value = 32004
$ws.Range("L83").Value = 4408.71435
$ws.Range("M83").Value = -27012
$ws.Range("N83").Value = -14392.71435
$ws.Range("H92").Value = 826.6316
$ws.Range("I92").Value = 753.2941
$ws.Range("K92").Value = 753.2941
$ws.Range("M92").Value = 494.7059
$ws.Range("H98").Value = 1592.1666
$ws.Range("I98").Value = 1139.3125
$ws.Range("J98").Value = 5215
$ws.Range("K98").Value = 1139.3125
$ws.Range("L98").Value = 5215
$ws.Range("M98").Value = 358.6875
$ws.Range("N98").Value = -8211
$ws.Range("H106").Value = 2715.8333
$ws.Range("I106").Value = 3027.2856
$ws.Range("J106").Value = 2279.8
$ws.Range("K106").Value = 3027.2856
$ws.Range("L106").Value = 2279.8
$ws.Range("M106").Value = -2396.2856
$ws.Range("N106").Value = -3541.8
$ws.Range("H107").Value = 599.7
$ws.Range("I107").Value = 652.1429000000001
$ws.Range("J107").Value = 477.33334
$ws.Range("K107").Value = 652.1429000000001
$ws.Range("L107").Value = 477.33334
$ws.Range("M107").Value = 1267.8571
$ws.Range("N107").Value = -4317.33334
$ws.Range("H116").Value = 4972.3887
$ws.Range("I116").Value = 4969.1875
$ws.Range("K116").Value = 4969.1875
$ws.Range("M116").Value = -1527.1875
$ws.Range("H122").Value = 1592.1666
$ws.Range("I122").Value = 1139.3125
$ws.Range("J122").Value = 5215
$ws.Range("K122").Value = 3417.9375
$ws.Range("L122").Value = 15645
$ws.Range("M122").Value = -967.9375
$ws.Range("N122").Value = -20545
$ws.Range("H137").Value = 2262.72
$ws.Range("I137").Value = 1889.5
$ws.Range("J137").Value = 4999.6665
$ws.Range("K137").Value = 5668.5
$ws.Range("L137").Value = 14998.9995
$ws.Range("M137").Value = -3118.5
$ws.Range("N137").Value = -20098.9995
$ws.Range("H139").Value = 104999
$ws.Range("J139").Value = 104999
$ws.Range("L139").Value = 104999
$ws.Range("N139").Value = -115279
$ws.Range("H140").Value = 99995
$ws.Range("J140").Value = 99995
$ws.Range("L140").Value = 99995
$ws.Range("N140").Value = -110355

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1478.1111
$ws.Range("I5").Value = 1006.625
$ws.Range("K5").Value = 1006.625
$ws.Range("M5").Value = -894.625
$ws.Range("H32").Value = 1822.2157
$ws.Range("I32").Value = 1785.6383
$ws.Range("K32").Value = 1785.6383
$ws.Range("M32").Value = -1498.6383
$ws.Range("H102").Value = 8335283
$ws.Range("I102").Value = 9092763
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 9092763
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -9091141
$ws.Range("N102").Value = -6244
$ws.Range("H132").Value = 4351159.5
$ws.Range("I132").Value = 5003168.5
$ws.Range("K132").Value = 15009505.5
$ws.Range("M132").Value = -15006975.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1478.1111
$ws.Range("I4").Value = 1006.625
$ws.Range("K4").Value = 1006.625
$ws.Range("M4").Value = -891.625
$ws.Range("H20").Value = 1800
$ws.Range("I20").Value = 1800
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1800
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -1553
$ws.Range("N20").ClearContents()
$ws.Range("H26").Value = 3471
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H86").Value = 3403.2173
$ws.Range("I86").Value = 3207.923
$ws.Range("J86").Value = 3657.1
$ws.Range("K86").Value = 3207.923
$ws.Range("L86").Value = 3657.1
$ws.Range("M86").Value = -2084.923
$ws.Range("N86").Value = -5903.1
$ws.Range("H89").Value = 3403.2173
$ws.Range("I89").Value = 3207.923
$ws.Range("J89").Value = 3657.1
$ws.Range("K89").Value = 16039.615
$ws.Range("L89").Value = 18285.5
$ws.Range("M89").Value = -10423.615
$ws.Range("N89").Value = -29517.5
$ws.Range("H100").Value = 15321.5
$ws.Range("J100").Value = 17762
$ws.Range("L100").Value = 17762
$ws.Range("N100").Value = -19926
$ws.Range("H103").Value = 22747.5
$ws.Range("J103").Value = 22747.5
$ws.Range("L103").Value = 22747.5
$ws.Range("N103").Value = -25091.5
$ws.Range("H134").Value = 20006726
$ws.Range("I134").Value = 20840152
$ws.Range("K134").Value = 62520456
$ws.Range("M134").Value = -62517921

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 6580.8887
$ws.Range("I7").Value = 9078.691999999999
$ws.Range("J7").Value = 86.59999999999999
$ws.Range("K7").Value = 9078.691999999999
$ws.Range("L7").Value = 86.59999999999999
$ws.Range("M7").Value = -8965.691999999999
$ws.Range("N7").Value = -312.6
$ws.Range("H31").Value = 5518.9443
$ws.Range("I31").Value = 3556.0667
$ws.Range("K31").Value = 3556.0667
$ws.Range("M31").Value = -3261.0667
$ws.Range("H34").Value = 5518.9443
$ws.Range("I34").Value = 3556.0667
$ws.Range("K34").Value = 3556.0667
$ws.Range("M34").Value = -3354.0667
$ws.Range("H93").Value = 21459.857
$ws.Range("I93").Value = 21459.857
$ws.Range("K93").Value = 21459.857
$ws.Range("M93").Value = -19587.857
$ws.Range("H95").Value = 13508
$ws.Range("J95").Value = 13508
$ws.Range("L95").Value = 13508
$ws.Range("N95").Value = -19000
$ws.Range("H96").Value = 16432.5
$ws.Range("J96").Value = 16432.5
$ws.Range("L96").Value = 16432.5
$ws.Range("N96").Value = -21924.5
$ws.Range("H111").Value = 29999
$ws.Range("J111").Value = 29999
$ws.Range("L111").Value = 29999
$ws.Range("N111").Value = -38179

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 49.5
$ws.Range("I2").Value = 45.5
$ws.Range("K2").Value = 273
$ws.Range("M2").Value = -160
$ws.Range("H25").Value = 10000
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H30").Value = 10000
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H121").Value = 245999.8
$ws.Range("I121").Value = 257499.75
$ws.Range("K121").Value = 772499.25
$ws.Range("M121").Value = -771189.25
$ws.Range("H122").Value = 1581.3636
$ws.Range("I122").Value = 1419.6
$ws.Range("J122").Value = 1716.1666
$ws.Range("K122").Value = 12776.4
$ws.Range("L122").Value = 15445.4994
$ws.Range("M122").Value = -10326.4
$ws.Range("N122").Value = -20345.4994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 28999.5
$ws.Range("I10").Value = 28999.5
$ws.Range("K10").Value = 28999.5
$ws.Range("M10").Value = -28830.5
$ws.Range("H62").Value = 32000
$ws.Range("I62").Value = 32000
$ws.Range("K62").Value = 32000
$ws.Range("M62").Value = -31314
$ws.Range("H65").Value = 32000
$ws.Range("I65").Value = 32000
$ws.Range("K65").Value = 96000
$ws.Range("M65").Value = -92568
$ws.Range("H70").Value = 6303.2856
$ws.Range("I70").Value = 6274.75
$ws.Range("K70").Value = 6274.75
$ws.Range("M70").Value = -6004.75
$ws.Range("H73").Value = 6303.2856
$ws.Range("I73").Value = 6274.75
$ws.Range("K73").Value = 6274.75
$ws.Range("M73").Value = -5338.75
$ws.Range("H104").Value = 29000
$ws.Range("J104").Value = 29000
$ws.Range("L104").Value = 29000
$ws.Range("N104").Value = -35988

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4129.8
$ws.Range("I7").Value = 4129.8
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 4129.8
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4017.8
$ws.Range("N7").ClearContents()
$ws.Range("H14").Value = 6999
$ws.Range("I14").Value = 1999
$ws.Range("J14").Value = 11999
$ws.Range("K14").Value = 1999
$ws.Range("L14").Value = 11999
$ws.Range("M14").Value = -1827
$ws.Range("N14").Value = -12343
$ws.Range("H40").Value = 3900.5334
$ws.Range("I40").Value = 3900.5334
$ws.Range("K40").Value = 3900.5334
$ws.Range("M40").Value = -3764.5334
$ws.Range("H126").Value = 4129.8
$ws.Range("I126").Value = 4129.8
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 12389.4
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -9919.400000000001
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 20837770
$ws.Range("I132").Value = 26318158
$ws.Range("J132").Value = 12300.4
$ws.Range("K132").Value = 78954474
$ws.Range("L132").Value = 36901.2
$ws.Range("M132").Value = -78951944
$ws.Range("N132").Value = -41961.2
